# "Generate Report for Handback" -- refresh the localization-status report
# after a successful handback: the zh-cn / de-de rows move from
# "Ready for handoff" to "Handed back: in sync with en-US", their handback
# timestamps advance, and the stale "handback file is not latest" error
# detail is cleared since the handback now matches the latest source.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status: handback completed, in sync with en-US ---
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime advances to the new handback run ---
$zhcn.Range("K2").Value = "2016-10-17 13:54:36"
$dede.Range("K2").Value = "2016-10-17 13:55:15"

# --- Error Detail cleared now that the handback is up to date ---
# (re-apply the cell's existing default style so the now-empty cell is
# still written out rather than dropped from the sheet)
$zhcnP2Style = $zhcn.Range("P2").Style
$zhcn.Range("P2").Value = ""
$zhcn.Range("P2").Style = $zhcnP2Style

$dedeP2Style = $dede.Range("P2").Style
$dede.Range("P2").Value = ""
$dede.Range("P2").Style = $dedeP2Style

# --- Column widths re-sized for the wider Status / Error Detail text ---
$overview.Columns.Item(5).ColumnWidth = 29.14437166849777
$overview.Columns.Item(6).ColumnWidth = 29.14437166849777

$zhcn.Columns.Item(3).ColumnWidth = 29.14437166849777
$zhcn.Columns.Item(16).ColumnWidth = 12.913719813028965

$dede.Columns.Item(3).ColumnWidth = 29.14437166849777
$dede.Columns.Item(16).ColumnWidth = 12.913719813028965
